$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns B..V (row 1)
$headers = @(
    "('N', 'C', 1)",
    "('S', 'H', 1)",
    "('C', 'S', 1.5)",
    "('O', 'C', 1.5)",
    "('C', 'O', 2)",
    "('C', 'O', 1)",
    "('C', 'C', 1)",
    "('C', 'C', 1.5)",
    "('C', 'H', 1)",
    "C",
    "H",
    "O",
    "amines",
    "ether",
    "aldehyde",
    "sulfide",
    "aromatic",
    "cycle",
    "weight",
    "logP",
    "SAS"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# New numeric data for rows 2..5, columns B..V
$row2 = @(3.507, 0.19, 0.368, 0.157, 0.874, 2.316, 0.678, 9.574999999999999, 12.021, 13.436, 13.599, 2.777, 1.407, 1.067, 0.874, 0.083, 1.933, 0.013, 187.431, 2.123, 2.319)
$row3 = @(7.642, 0.05, 0.095, 0.148, 1.317, 0.456, 18.148, 1.216, 34.955, 21.98, 36.906, 1.471, 2.51, 0.187, 1.317, 0.029, 0.367, 3.418, 280.898, 2.338, 4.402)
$row4 = @(0.644, 0.036, 0.074, 0.19, 0.6, 7.832, 7.493, 6.12, 26.169, 16.944, 26.64, 5.717, 0.338, 3.646, 0.6, 0.034, 1.02, 1.018, 225.663, 1.889, 3.242)
$row5 = @(13.453, 0.021, 0.018, 0.028, 2.195, 0.635, 8.119, 1.113, 22.878, 15.945, 25.361, 2.347, 4.334, 0.374, 2.195, 0.024, 0.08799999999999999, 2.725, 266.444, 0.472, 3.101)

$rows = @($row2, $row3, $row4, $row5)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowValues = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 2).Value = $rowValues[$c]
    }
}

# Remove now-unused trailing columns W and X (old sheet had 2 more columns than the new one)
$ws.Columns("W:X").Delete()
